# This script reproduces the commit:
#   "Inclui testes para campos de ano da norma, data de promulgação,
#   título da norma e seção do DOU. Pequenos ajustes nas funções de
#   limpar texto e separar portarias multiplas"
#
# Net effect on the worksheet (gabarito.xlsx, sheet "Gabarito"):
#   - A handful of NR_ATO/título text cells (column I, mostly) get their
#     text "cleaned up" (trailing period added, ordinal "º" added after a
#     lone "1" day-of-month, thousands separator "." added to 4-digit
#     portaria numbers).
#   - The sheet view scrolled down / re-selected a different cell.
#
# The edits below are applied in the same order the cells' text was
# actually changed upstream, which matters because each newly-introduced
# (not-previously-seen) string value is appended to the workbook's
# shared-string table in edit order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clean up / fix individual título (NR_ATO display) strings ---------

# Trailing "." added to a blank Instrução Normativa entry.
$ws.Range("I7").Value = "INSTRUÇÃO NORMATIVA Nº , DE DE 2017."

# Add missing "." thousands separator to four-digit portaria numbers.
$ws.Range("I70").Value = "PORTARIA Nº 1.931, DE 23 DE AGOSTO DE 2016"
$ws.Range("I69").Value = "PORTARIA Nº 1.930, DE 23 DE AGOSTO DE 2016"
$ws.Range("I68").Value = "PORTARIA Nº 1.927, DE 23 DE AGOSTO DE 2016"

# Add ordinal marker "º" after "DE 1" (1st of June) for this batch of
# portarias published multiple-at-once on 2017-06-01.
$ws.Range("I31").Value = "PORTARIA Nº 1.242, DE 1º DE JUNHO DE 2017"
$ws.Range("I35").Value = "PORTARIA Nº 1.245, DE 1º DE JUNHO DE 2017"
$ws.Range("I36").Value = "PORTARIA Nº 1.246, DE 1º DE JUNHO DE 2017"
$ws.Range("I37").Value = "PORTARIA Nº 1.256, DE 1º DE JUNHO DE 2017"
$ws.Range("I32").Value = "PORTARIA Nº 1.252, DE 1º DE JUNHO DE 2017"
$ws.Range("I33").Value = "PORTARIA Nº 1.253, DE 1º DE JUNHO DE 2017"
$ws.Range("I34").Value = "PORTARIA Nº 1.254, DE 1º DE JUNHO DE 2017"

# Trailing "." added.
$ws.Range("I61").Value = "PORTARIA Nº 158, DE 24 DE AGOSTO DE 2016."

# --- Scroll the view down and move the selection ------------------------

$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("I62").Select()
